$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.183.77'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.352.51'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.45'
$ws.Range("D5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.62'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.96%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.89%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.12'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.66%  '

# Row 11
$ws.Range("E11").Value = '  +0.10%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.72'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.99%  '

# Row 13
$ws.Range("E13").Value = '  +3.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.73'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.77%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.721.73'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.55%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.338.96'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.60%  '

# Row 17
$ws.Range("E17").Value = '  +1.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.164.97'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.18'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.37%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.25'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.14%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.17'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.50'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.57%  '

# Row 25
$ws.Range("E25").Value = '  -0.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.45%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.53'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.11%  '

# Row 28
$ws.Range("E28").Value = '  +7.10%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.58%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.40'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.22%  '

# Row 31
$ws.Range("E31").Value = '  -0.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.02'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.20%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0724'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.41%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.16'
$ws.Range("D34").ClearFormats()

# Row 35
$ws.Range("E35").Value = '  -0.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.71%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.32'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.59%  '

# Row 38
$ws.Range("E38").Value = '  +0.37%  '

# Row 39
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +13.92%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.73%  '

# Row 41
$ws.Range("E41").Value = '  -0.10%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.02'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -33.15%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.942.91'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.42%  '

# Row 44
$ws.Range("E44").Value = '  +0.29%  '

# Row 45
$ws.Range("E45").Value = '  +3.62%  '

# Row 46
$ws.Range("E46").Value = '  -10.07%  '

# Row 47
$ws.Range("E47").Value = '  -0.76%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.585.14'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.90'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.33%  '

# Row 50
$ws.Range("E50").Value = '  -1.88%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.06'
$ws.Range("D51").ClearFormats()
